# Update market-price/profit figures pulled by the scheduled pricing runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1087429.1
$ws.Range("J17").Value = 1087429.1
$ws.Range("L17").Value = 3262287.3
$ws.Range("N17").Value = -3262623.3
$ws.Range("H76").Value = 3008.0168
$ws.Range("I76").Value = 2667.7058
$ws.Range("K76").Value = 2667.7058
$ws.Range("M76").Value = -2352.7058
$ws.Range("H79").Value = 3008.0168
$ws.Range("I79").Value = 2667.7058
$ws.Range("K79").Value = 2667.7058
$ws.Range("M79").Value = -1575.7058
$ws.Range("H93").Value = 88957.664
$ws.Range("J93").Value = 88957.664
$ws.Range("L93").Value = 88957.664
$ws.Range("N93").Value = -93949.664
$ws.Range("H103").Value = 477.33334
$ws.Range("I103").Value = 400
$ws.Range("J103").Value = 516
$ws.Range("K103").Value = 1200
$ws.Range("L103").Value = 1548
$ws.Range("M103").Value = -614
$ws.Range("N103").Value = -2720
$ws.Range("H112").Value = 1203.2258
$ws.Range("I112").Value = 1400
$ws.Range("J112").Value = 1189.6552
$ws.Range("K112").Value = 4200
$ws.Range("L112").Value = 3568.9656
$ws.Range("M112").Value = -3092
$ws.Range("N112").Value = -5784.9656
$ws.Range("H137").Value = 1576.174
$ws.Range("I137").Value = 1640.6428
$ws.Range("J137").Value = 1475.8889
$ws.Range("K137").Value = 4921.928400000001
$ws.Range("L137").Value = 4427.6667
$ws.Range("M137").Value = -2371.928400000001
$ws.Range("N137").Value = -9527.6667
$ws.Range("H138").Value = 2677.9326
$ws.Range("I138").Value = 1342.1754
$ws.Range("J138").Value = 5057.25
$ws.Range("K138").Value = 4026.5262
$ws.Range("L138").Value = 15171.75
$ws.Range("M138").Value = 1113.4738
$ws.Range("N138").Value = -25451.75
$ws.Range("H141").Value = 5210.052
$ws.Range("I141").Value = 1252.8043
$ws.Range("J141").Value = 20379.5
$ws.Range("K141").Value = 3758.4129
$ws.Range("L141").Value = 61138.5
$ws.Range("M141").Value = 1421.5871
$ws.Range("N141").Value = -71498.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5093.482
$ws.Range("I32").Value = 3597
$ws.Range("K32").Value = 3597
$ws.Range("M32").Value = -3310
$ws.Range("H45").Value = 1304.875
$ws.Range("I45").Value = 1075
$ws.Range("J45").Value = 1994.5
$ws.Range("K45").Value = 1075
$ws.Range("L45").Value = 1994.5
$ws.Range("M45").Value = -698
$ws.Range("N45").Value = -2748.5
$ws.Range("H74").Value = 3811.561
$ws.Range("I74").Value = 4061.081
$ws.Range("J74").Value = 1503.5
$ws.Range("K74").Value = 4061.081
$ws.Range("L74").Value = 1503.5
$ws.Range("M74").Value = -3187.081
$ws.Range("N74").Value = -3251.5
$ws.Range("H77").Value = 3811.561
$ws.Range("I77").Value = 4061.081
$ws.Range("J77").Value = 1503.5
$ws.Range("K77").Value = 20305.405
$ws.Range("L77").Value = 7517.5
$ws.Range("M77").Value = -15937.405
$ws.Range("N77").Value = -16253.5
$ws.Range("H97").Value = 1423.76
$ws.Range("I97").Value = 932.8889
$ws.Range("J97").Value = 2686
$ws.Range("K97").Value = 932.8889
$ws.Range("L97").Value = 2686
$ws.Range("M97").Value = -436.8889
$ws.Range("N97").Value = -3678
$ws.Range("H122").Value = 1514.1818
$ws.Range("I122").Value = 1167.9166
$ws.Range("J122").Value = 2437.5557
$ws.Range("K122").Value = 3503.7498
$ws.Range("L122").Value = 7312.6671
$ws.Range("M122").Value = -1053.7498
$ws.Range("N122").Value = -12212.6671
$ws.Range("H132").Value = 1973.6666
$ws.Range("I132").Value = 1008.73914
$ws.Range("J132").Value = 2982.4546
$ws.Range("K132").Value = 3026.21742
$ws.Range("L132").Value = 8947.363799999999
$ws.Range("M132").Value = -496.2174199999999
$ws.Range("N132").Value = -14007.3638

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 800
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 800
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 800
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -1250
$ws.Range("H67").Value = 800
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 800
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 800
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -2360
$ws.Range("H134").Value = 1396.3158
$ws.Range("I134").Value = 1087.742
$ws.Range("J134").Value = 2762.8572
$ws.Range("K134").Value = 3263.226
$ws.Range("L134").Value = 8288.571599999999
$ws.Range("M134").Value = -728.2259999999997
$ws.Range("N134").Value = -13358.5716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2577.677
$ws.Range("I31").Value = 1641.125
$ws.Range("J31").Value = 3485.8484
$ws.Range("K31").Value = 1641.125
$ws.Range("L31").Value = 3485.8484
$ws.Range("M31").Value = -1346.125
$ws.Range("N31").Value = -4075.8484
$ws.Range("H34").Value = 2577.677
$ws.Range("I34").Value = 1641.125
$ws.Range("J34").Value = 3485.8484
$ws.Range("K34").Value = 1641.125
$ws.Range("L34").Value = 3485.8484
$ws.Range("M34").Value = -1439.125
$ws.Range("N34").Value = -3889.8484
$ws.Range("H93").Value = 19289.445
$ws.Range("I93").Value = 1841
$ws.Range("K93").Value = 1841
$ws.Range("M93").Value = 31
$ws.Range("H99").Value = 11993.272
$ws.Range("I99").Value = 2532.4
$ws.Range("J99").Value = 19877.334
$ws.Range("K99").Value = 2532.4
$ws.Range("L99").Value = 19877.334
$ws.Range("M99").Value = -1034.4
$ws.Range("N99").Value = -22873.334
$ws.Range("H126").Value = 11993.272
$ws.Range("I126").Value = 2532.4
$ws.Range("J126").Value = 19877.334
$ws.Range("K126").Value = 7597.200000000001
$ws.Range("L126").Value = 59632.00199999999
$ws.Range("M126").Value = -5127.200000000001
$ws.Range("N126").Value = -64572.00199999999
$ws.Range("H132").Value = 1985.2903
$ws.Range("I132").Value = 1009.1875
$ws.Range("J132").Value = 3026.4666
$ws.Range("K132").Value = 3027.5625
$ws.Range("L132").Value = 9079.399800000001
$ws.Range("M132").Value = -497.5625
$ws.Range("N132").Value = -14139.3998
$ws.Range("H134").Value = 1236.8889
$ws.Range("I134").Value = 1104.8197
$ws.Range("J134").Value = 1969.2727
$ws.Range("K134").Value = 3314.4591
$ws.Range("L134").Value = 5907.8181
$ws.Range("M134").Value = -779.4591
$ws.Range("N134").Value = -10977.8181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1412.5714
$ws.Range("I97").Value = 1585.2727
$ws.Range("J97").Value = 779.3333
$ws.Range("K97").Value = 1585.2727
$ws.Range("L97").Value = 779.3333
$ws.Range("M97").Value = -1089.2727
$ws.Range("N97").Value = -1771.3333
$ws.Range("H102").Value = 3882.4
$ws.Range("I102").Value = 912
$ws.Range("J102").Value = 4625
$ws.Range("K102").Value = 912
$ws.Range("L102").Value = 4625
$ws.Range("M102").Value = 710
$ws.Range("N102").Value = -7869
$ws.Range("H132").Value = 1421.871
$ws.Range("I132").Value = 1421.871
$ws.Range("K132").Value = 4265.613
$ws.Range("M132").Value = -1735.613

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 23764.08
$ws.Range("I40").Value = 29179.053
$ws.Range("J40").Value = 6616.6665
$ws.Range("K40").Value = 29179.053
$ws.Range("L40").Value = 6616.6665
$ws.Range("M40").Value = -29043.053
$ws.Range("N40").Value = -6888.6665
$ws.Range("H61").Value = 2584
$ws.Range("I61").Value = 2376
$ws.Range("K61").Value = 2376
$ws.Range("M61").Value = -2174
$ws.Range("H100").Value = 1426.3334
$ws.Range("I100").Value = 779
$ws.Range("K100").Value = 779
$ws.Range("M100").Value = -238
$ws.Range("H113").Value = 2584
$ws.Range("I113").Value = 2376
$ws.Range("K113").Value = 2376
$ws.Range("M113").Value = -206
$ws.Range("H136").Value = 12822602
$ws.Range("I136").Value = 2243.9583
$ws.Range("K136").Value = 6731.874899999999
$ws.Range("M136").Value = -4181.874899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 870.88
$ws.Range("I113").Value = 908.1579
$ws.Range("J113").Value = 752.8333
$ws.Range("K113").Value = 2724.4737
$ws.Range("L113").Value = 2258.4999
$ws.Range("M113").Value = -554.4737
$ws.Range("N113").Value = -6598.4999
$ws.Range("H122").Value = 64759.625
$ws.Range("I122").Value = 102066
$ws.Range("J122").Value = 2582.3333
$ws.Range("K122").Value = 306198
$ws.Range("L122").Value = 7746.999899999999
$ws.Range("M122").Value = -303748
$ws.Range("N122").Value = -12646.9999
$ws.Range("H132").Value = 1461.7167
$ws.Range("I132").Value = 1282.0952
$ws.Range("K132").Value = 3846.2856
$ws.Range("M132").Value = -1316.2856
$ws.Range("H136").Value = 6291048
$ws.Range("I136").Value = 9524444
$ws.Range("K136").Value = 28573332
$ws.Range("M136").Value = -28570782
